# Correccion cola de prioridad
# D13 currently holds the text "21" (shared string index 7); it should be
# the text "15" (the priority/age value that the queue entry actually
# corresponds to). The cell must stay a plain text cell (no numeric
# auto-conversion, no style change) to match the original authoring.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$cell = $ws.Range("D13")

# Force the value in as text (leading apostrophe keeps Excel from
# reinterpreting "15" as a number) ...
$cell.Value = "'15"

# ... then drop the quote-prefix / number-format bookkeeping that the
# apostrophe entry leaves behind, so the cell's style stays identical to
# before (General, no explicit style index).
$cell.ClearFormats()
